$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51.
# D-column values are plain text (e.g. "24.603.37", "1.003") that Excel would
# otherwise auto-coerce into numbers; force text storage via NumberFormat="@"
# then restore the default "Normal" style so no stray number-format survives.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "24.603.37"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.691.90"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.51"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3896"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.40%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4041"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.70%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "1.495"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.04%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "52.89"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08755"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.71%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "25.58"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +7.95%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.523"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.43%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.00001356"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.52%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.961"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.51%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.689.54"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "98.59"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.07109"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.07%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "19.97"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.295"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.01%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "14.29"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "24.602.35"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.002"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -6.53%  "

$ws.Range("E26").Value = "  -0.42%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "22.80"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "162.25"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.758"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +15.31%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "137.26"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("E31").Value = "  +0.63%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.873.75"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.08849"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.52%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.411"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.16%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.037"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.07%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.969"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.38%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02937"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +7.60%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2752"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.56%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "10.81"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.86%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "14.29"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.68%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.09140"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.85%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.7905"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.79%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.462"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.75"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.94%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.7219"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.578"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.209"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("E48").Value = "  -0.07%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.342"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.39%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "137.84"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "91.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
